$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report date range) ---
$ws.Range("A8").Value = 'Volume 30   Number  22'
$ws.Range("C9").Value = 'Report Covering the Week  5/29/2023  Through  6/4/2023'

# --- Crime statistics table updates (rows 16-27) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 5
$ws.Range("L16").Value = -16.666666666666
$ws.Range("M16").Value = -44.444444444444
$ws.Range("N16").Value = -91.228070175438
$ws.Range("J17").Value = 6
$ws.Range("K17").Value = -50
$ws.Range("N17").Value = -84.210526315789
$ws.Range("D18").Value = '''0'
$ws.Range("E18").Value = '***.*'
$ws.Range("F18").Value = 2
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 3
$ws.Range("K18").Value = 50
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = -76.923076923076
$ws.Range("C19").Value = 2
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 50
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 7
$ws.Range("L19").Value = 16.666666666666
$ws.Range("M19").Value = -44
$ws.Range("N19").Value = -75.862068965517
$ws.Range("C21").Value = 3
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 9
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 12.5
$ws.Range("I21").Value = 25
$ws.Range("J21").Value = 21
$ws.Range("K21").Value = 19.047619047619
$ws.Range("L21").Value = 8.695652173913
$ws.Range("M21").Value = -32.432432432432
$ws.Range("N21").Value = -83.660130718954
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = '''0'
$ws.Range("E24").Value = '***.*'
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 14
$ws.Range("J24").Value = 9
$ws.Range("K24").Value = 55.555555555555
$ws.Range("L24").Value = 7.692307692307
$ws.Range("M24").Value = -46.153846153846
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 60
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 13
$ws.Range("K25").Value = 23.076923076923
$ws.Range("L25").Value = 128.571428571429
$ws.Range("M25").Value = 33.333333333333
$ws.Range("C27").Value = '''0'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = '''0'
$ws.Range("H27").Value = '***.*'
